$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("fragecodes")
$ws1.Range("B28").Value = "time"
$ws1.Range("B29").Value = "time"
$ws1.Range("B35").Value = "time"
$ws1.Range("B36").Value = "time"

$ws3 = $wb.Worksheets.Item("invites")
$ws3.Range("B49").Value = 80000
$ws3.Range("B49").Interior.Color = 65535
$ws3.Activate()
[void]$ws3.Range("B49").Select()

[void]$ws1.Activate()
[void]$ws1.Range("C37").Select()
